# "Clean commit after removing API key" — refresh the product dummy-data
# color palette (drop pink/purple/orange, add mustard/dark green/navy/
# "black " as new shared strings) and re-point each product row's
# size (P) / color (Q) cells at the refreshed palette. Also nudges the
# saved view/selection and widens column A to fit the longer color names.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New size (column P) values for rows 2-26 (unchanged semantically —
# written back explicitly so the shared-string table is rebuilt cleanly)
$sizes = @{
    2  = "small";  3  = "small";  4  = "small"
    5  = "large";  6  = "large";  7  = "large";  8  = "large"
    9  = "medium"; 10 = "medium"; 11 = "medium"; 12 = "medium"; 13 = "medium"
    14 = "small"
    15 = "xlarge"; 16 = "xlarge"; 17 = "xlarge"
    18 = "xxlarge"
    19 = "small";  20 = "small";  21 = "small"
    22 = "medium"; 23 = "medium"; 24 = "medium"; 25 = "medium"; 26 = "medium"
}

# New color (column Q) values for rows 2-26 — pink/purple/orange retired,
# replaced across the sheet with mustard / dark green / navy / "black "
$colors = @{
    2  = "mustard";    3  = "mustard";    4  = "white"
    5  = "blue";       6  = "green";      7  = "white";      8  = "black"
    9  = "dark green"; 10 = "white";      11 = "blue";       12 = "yellow"; 13 = "yellow"
    14 = "white"
    15 = "yellow";     16 = "blue";       17 = "green"
    18 = "red"
    19 = "navy";       20 = "yellow";     21 = "red"
    22 = "mustard";    23 = "black ";     24 = "green";      25 = "green";  26 = "white"
}

foreach ($row in 2..26) {
    $ws.Range("P$row").Value = $sizes[$row]
    $ws.Range("Q$row").Value = $colors[$row]
}

# Widen column A for the new, longer color labels
$ws.Columns.Item(1).ColumnWidth = 39.19921875

# Restore the saved selection/scroll position
$ws.Range("T19").Select()
